# Applies the "Atualizado por script em 12-11-2023 02:45" update:
#   1) Swap the match-detail columns (F:V) between row pairs 62/63, 66/67,
#      104/105 and 133/134 (the underlying fixture list got re-sorted by
#      kickoff time; columns A:E - index/pais/torneio/temporada/data - stay
#      put since the date grouping did not change).
#   2) Append three newly scraped fixtures as rows 173-175.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($r1, $r2) {
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Swap-Row 62 63
Swap-Row 66 67
Swap-Row 104 105
Swap-Row 133 134

# --- Append three new fixtures as rows 173-175. ---------------------------
# Write the values first (column D - "temporada" - is forced to Text so the
# "2023" string isn't silently re-typed as a number), then copy the visual
# style (bold/centered "Indice" column, datetime format on "data_partida")
# from the last existing row so the new rows look like the rest of the table.

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Set-MatchRow($r, $values) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"   # keep "temporada" as text, e.g. "2023"
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $values[$i]
    }
}

Set-MatchRow 173 @(
    172, "argentina", "copa-de-la-liga-profesional", "2023", 45241.9375,
    "Rosario Central", 3, "River Plate", 1,
    3.49, "08/11/2023 14:42", 4.63, "11/11/2023 22:28",
    3.28, "08/11/2023 14:42", 3.89, "11/11/2023 22:28",
    2.13, "08/11/2023 14:42", 1.78, "11/11/2023 22:24",
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/rosario-central-river-plate/SnhHrHYJ/"
)

Set-MatchRow 174 @(
    173, "argentina", "copa-de-la-liga-profesional", "2023", 45242,
    "Lanus", 0, "Racing Club", 2,
    2.79, "08/11/2023 14:42", 2.48, "11/11/2023 23:59",
    3.15, "08/11/2023 14:42", 3.11, "11/11/2023 23:59",
    2.61, "08/11/2023 14:42", 3.24, "11/11/2023 23:59",
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/lanus-racing-club/UyGku3qP/"
)

Set-MatchRow 175 @(
    174, "argentina", "copa-de-la-liga-profesional", "2023", 45242.04166666666,
    "Huracan", 1, "Arsenal Sarandi", 0,
    1.64, "08/11/2023 14:42", 1.66, "12/11/2023 00:58",
    3.58, "08/11/2023 14:42", 3.37, "12/11/2023 00:59",
    5.67, "08/11/2023 14:42", 7.26, "12/11/2023 00:58",
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/huracan-arsenal-sarandi/6wiLsylQ/"
)

# Re-apply the row-172 look (bold/bordered index column, datetime format)
# to the newly-written rows without disturbing the values just written.
$ws.Range("A172:V172").Copy()
$ws.Range("A173:V175").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A1").Select()
